# Apply the "IBM -> LinuxForHealth" rebrand edit to the
# StructureDefinition-reengagement-display workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet: URL, Version, Date, Publisher
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reengagement-display"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------
# Elements sheet: Fixed Value for Extension.url (row 5), and clear the
# stray duplicated constraint text that had been listed against the
# root Extension row (row 2) - it correctly belongs only on the
# Extension.extension row.
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reengagement-display"
$elements.Range("AI2").Value = ""
